# Scheduled-runner style update: refresh market-board price figures and
# derived profit columns (H..N) across the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 309.66666
$ws.Range("J2").Value = 234
$ws.Range("L2").Value = 234
$ws.Range("N2").Value = -460
$ws.Range("H132").Value = 16433.857
$ws.Range("J132").Value = 700
$ws.Range("L132").Value = 2100
$ws.Range("N132").Value = -7160
$ws.Range("H137").Value = 2714.3584
$ws.Range("I137").Value = 1810
$ws.Range("K137").Value = 5430
$ws.Range("M137").Value = -2880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4127.4287
$ws.Range("I61").Value = 4127.4287
$ws.Range("K61").Value = 4127.4287
$ws.Range("M61").Value = -3915.4287
$ws.Range("H136").Value = 4127.4287
$ws.Range("I136").Value = 4127.4287
$ws.Range("K136").Value = 12382.2861
$ws.Range("M136").Value = -9832.286100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2095
$ws.Range("I20").Value = 2095
$ws.Range("K20").Value = 2095
$ws.Range("M20").Value = -1848
$ws.Range("H87").Value = 86666.664
$ws.Range("J87").Value = 86666.664
$ws.Range("L87").Value = 86666.664
$ws.Range("N87").Value = -89162.664
$ws.Range("H90").Value = 86666.664
$ws.Range("J90").Value = 86666.664
$ws.Range("L90").Value = 259999.992
$ws.Range("N90").Value = -272479.992
$ws.Range("H96").Value = 39999.5
$ws.Range("I96").Value = 39999.5
$ws.Range("K96").Value = 39999.5
$ws.Range("M96").Value = -37253.5
$ws.Range("H100").Value = 8546.666999999999
$ws.Range("J100").Value = 8546.666999999999
$ws.Range("L100").Value = 8546.666999999999
$ws.Range("N100").Value = -10710.667
$ws.Range("H134").Value = 1670.1428
$ws.Range("I134").Value = 1670.1428
$ws.Range("K134").Value = 5010.428400000001
$ws.Range("M134").Value = -2475.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2250
$ws.Range("I3").Value = 2250
$ws.Range("K3").Value = 2250
$ws.Range("M3").Value = -2137
$ws.Range("H6").Value = 1325
$ws.Range("I6").Value = 1325
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1325
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1212
$ws.Range("N6").ClearContents()
$ws.Range("H7").Value = 3909.3794
$ws.Range("J7").Value = 608.6667
$ws.Range("L7").Value = 608.6667
$ws.Range("N7").Value = -834.6667
$ws.Range("H16").Value = 625
$ws.Range("I16").Value = 605
$ws.Range("J16").Value = 665
$ws.Range("K16").Value = 605
$ws.Range("L16").Value = 665
$ws.Range("M16").Value = -318
$ws.Range("N16").Value = -1239
$ws.Range("H28").Value = 11549.333
$ws.Range("J28").Value = 11549.333
$ws.Range("L28").Value = 11549.333
$ws.Range("N28").Value = -12039.333
$ws.Range("H43").Value = 18059.8
$ws.Range("J43").Value = 18059.8
$ws.Range("L43").Value = 18059.8
$ws.Range("N43").Value = -18427.8
$ws.Range("H101").Value = 18059.8
$ws.Range("J101").Value = 18059.8
$ws.Range("L101").Value = 18059.8
$ws.Range("N101").Value = -24549.8
$ws.Range("H113").Value = 625
$ws.Range("I113").Value = 605
$ws.Range("J113").Value = 665
$ws.Range("K113").Value = 605
$ws.Range("L113").Value = 665
$ws.Range("M113").Value = 1565
$ws.Range("N113").Value = -5005
$ws.Range("H122").Value = 1341
$ws.Range("I122").Value = 1175.4166
$ws.Range("K122").Value = 3526.2498
$ws.Range("M122").Value = -1076.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3000
$ws.Range("J9").Value = 3000
$ws.Range("L9").Value = 9000
$ws.Range("N9").Value = -9448
$ws.Range("H14").Value = 767.4
$ws.Range("I14").Value = 767.4
$ws.Range("K14").Value = 2302.2
$ws.Range("M14").Value = -2129.2
$ws.Range("H34").Value = 972.8182
$ws.Range("I34").Value = 234
$ws.Range("J34").Value = 1249.875
$ws.Range("K34").Value = 702
$ws.Range("L34").Value = 3749.625
$ws.Range("M34").Value = -618
$ws.Range("N34").Value = -3917.625
$ws.Range("H39").Value = 6342.231
$ws.Range("J39").Value = 7438.091
$ws.Range("L39").Value = 22314.273
$ws.Range("N39").Value = -22902.273
$ws.Range("H109").Value = 1230.1428
$ws.Range("I109").Value = 1185.1666
$ws.Range("J109").Value = 1500
$ws.Range("K109").Value = 3555.4998
$ws.Range("L109").Value = 4500
$ws.Range("M109").Value = -2515.4998
$ws.Range("N109").Value = -6580
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H132").Value = 2471.2144
$ws.Range("J132").Value = 3083.1667
$ws.Range("L132").Value = 27748.5003
$ws.Range("N132").Value = -32808.5003
$ws.Range("H136").Value = 9666.666999999999
$ws.Range("I136").Value = 6000
$ws.Range("K136").Value = 18000
$ws.Range("M136").Value = -12900
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H92").Value = 5000
$ws.Range("J92").Value = 5000
$ws.Range("L92").Value = 5000
$ws.Range("N92").Value = -8744
$ws.Range("H97").Value = 744
$ws.Range("I97").Value = 573.3333
$ws.Range("K97").Value = 573.3333
$ws.Range("M97").Value = -77.33330000000001
$ws.Range("H102").Value = 2289.4614
$ws.Range("I102").Value = 2073.9092
$ws.Range("K102").Value = 2073.9092
$ws.Range("M102").Value = -451.9092000000001
$ws.Range("H113").Value = 7000
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 9333.333000000001
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 9333.333000000001
$ws.Range("M113").Value = -1330
$ws.Range("N113").Value = -13673.333
$ws.Range("H122").Value = 2582.4
$ws.Range("I122").Value = 2228
$ws.Range("K122").Value = 6684
$ws.Range("M122").Value = -4234
$ws.Range("H132").Value = 8499.5
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1924.7142
$ws.Range("I22").Value = 1759.8889
$ws.Range("J22").Value = 2221.4
$ws.Range("K22").Value = 1759.8889
$ws.Range("L22").Value = 2221.4
$ws.Range("M22").Value = -1464.8889
$ws.Range("N22").Value = -2811.4
$ws.Range("H27").Value = 1924.7142
$ws.Range("I27").Value = 1759.8889
$ws.Range("J27").Value = 2221.4
$ws.Range("K27").Value = 1759.8889
$ws.Range("L27").Value = 2221.4
$ws.Range("M27").Value = -1652.8889
$ws.Range("N27").Value = -2435.4
$ws.Range("H40").Value = 3196.6667
$ws.Range("I40").Value = 3845
$ws.Range("K40").Value = 3845
$ws.Range("M40").Value = -3709
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H63").Value = 67222
$ws.Range("I63").Value = 67222
$ws.Range("K63").Value = 67222
$ws.Range("M63").Value = -66473
$ws.Range("H66").Value = 67222
$ws.Range("I66").Value = 67222
$ws.Range("K66").Value = 201666
$ws.Range("M66").Value = -197922
$ws.Range("H122").Value = 3649.3333
$ws.Range("I122").Value = 3649.3333
$ws.Range("K122").Value = 10947.9999
$ws.Range("M122").Value = -8497.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 41799.6
$ws.Range("I54").Value = 3000
$ws.Range("K54").Value = 3000
$ws.Range("M54").Value = -2480
$ws.Range("H100").Value = 672.0714
$ws.Range("I100").Value = 610.9
$ws.Range("K100").Value = 1221.8
$ws.Range("M100").Value = -680.8
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 2800
$ws.Range("K122").Value = 8400
$ws.Range("M122").Value = -5950
